# Scheduled market-data refresh: update cached price/profit columns
# (H..N) for the affected leve rows across the ALC/ARM/CRP/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5089.6665
$ws.Range("I64").Value = 5169.1055
$ws.Range("J64").Value = 4901
$ws.Range("K64").Value = 5169.1055
$ws.Range("L64").Value = 4901
$ws.Range("M64").Value = -4921.1055
$ws.Range("N64").Value = -5397
$ws.Range("H67").Value = 5089.6665
$ws.Range("I67").Value = 5169.1055
$ws.Range("J67").Value = 4901
$ws.Range("K67").Value = 5169.1055
$ws.Range("L67").Value = 4901
$ws.Range("M67").Value = -4311.1055
$ws.Range("N67").Value = -6617
$ws.Range("H69").Value = 5816.1816
$ws.Range("I69").Value = 5816.1816
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 17448.5448
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -16574.5448
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 5816.1816
$ws.Range("I72").Value = 5816.1816
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 52345.6344
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -47977.6344
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 5762.7334
$ws.Range("I74").Value = 5499.9165
$ws.Range("K74").Value = 5499.9165
$ws.Range("M74").Value = -4563.9165
$ws.Range("H76").Value = 4360.3125
$ws.Range("I76").Value = 4365.4287
$ws.Range("K76").Value = 4365.4287
$ws.Range("M76").Value = -4050.4287
$ws.Range("H77").Value = 5762.7334
$ws.Range("I77").Value = 5499.9165
$ws.Range("K77").Value = 27499.5825
$ws.Range("M77").Value = -22819.5825
$ws.Range("H79").Value = 4360.3125
$ws.Range("I79").Value = 4365.4287
$ws.Range("K79").Value = 4365.4287
$ws.Range("M79").Value = -3273.4287
$ws.Range("H88").Value = 2768.4443
$ws.Range("I88").Value = 1990
$ws.Range("J88").Value = 2865.75
$ws.Range("K88").Value = 1990
$ws.Range("L88").Value = 2865.75
$ws.Range("M88").Value = -1584
$ws.Range("N88").Value = -3677.75
$ws.Range("H91").Value = 2768.4443
$ws.Range("I91").Value = 1990
$ws.Range("J91").Value = 2865.75
$ws.Range("K91").Value = 1990
$ws.Range("L91").Value = 2865.75
$ws.Range("M91").Value = -586
$ws.Range("N91").Value = -5673.75
$ws.Range("H99").Value = 611.2857
$ws.Range("I99").Value = 312.15384
$ws.Range("K99").Value = 936.4615200000001
$ws.Range("M99").Value = 561.5384799999999
$ws.Range("H101").Value = 526.5714
$ws.Range("I101").Value = 284.33334
$ws.Range("J101").Value = 1980
$ws.Range("K101").Value = 853.0000200000001
$ws.Range("L101").Value = 5940
$ws.Range("M101").Value = 768.9999799999999
$ws.Range("N101").Value = -9184
$ws.Range("H108").Value = 90342
$ws.Range("J108").Value = 90342
$ws.Range("L108").Value = 90342
$ws.Range("N108").Value = -98022
$ws.Range("H116").Value = 78197.39
$ws.Range("I116").Value = 88811.164
$ws.Range("K116").Value = 88811.164
$ws.Range("M116").Value = -85369.164
$ws.Range("H123").Value = 89995
$ws.Range("J123").Value = 89995
$ws.Range("L123").Value = 89995
$ws.Range("N123").Value = -99795
$ws.Range("H137").Value = 3916
$ws.Range("J137").Value = 5868.25
$ws.Range("L137").Value = 17604.75
$ws.Range("N137").Value = -22704.75
$ws.Range("H138").Value = 3180.7693
$ws.Range("I138").Value = 1022.53845
$ws.Range("J138").Value = 3900.1794
$ws.Range("K138").Value = 3067.61535
$ws.Range("L138").Value = 11700.5382
$ws.Range("M138").Value = 2072.38465
$ws.Range("N138").Value = -21980.5382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1267.375
$ws.Range("I45").Value = 1037.4
$ws.Range("J45").Value = 1650.6666
$ws.Range("K45").Value = 1037.4
$ws.Range("L45").Value = 1650.6666
$ws.Range("M45").Value = -660.4000000000001
$ws.Range("N45").Value = -2404.6666
$ws.Range("H76").Value = 45095.668
$ws.Range("J76").Value = 45095.668
$ws.Range("L76").Value = 45095.668
$ws.Range("N76").Value = -45771.668
$ws.Range("H79").Value = 45095.668
$ws.Range("J79").Value = 45095.668
$ws.Range("L79").Value = 45095.668
$ws.Range("N79").Value = -47435.668
$ws.Range("H92").Value = 15150
$ws.Range("J92").Value = 15150
$ws.Range("L92").Value = 15150
$ws.Range("N92").Value = -20142
$ws.Range("H122").Value = 2550.5908
$ws.Range("I122").Value = 2426.8
$ws.Range("K122").Value = 7280.400000000001
$ws.Range("M122").Value = -4830.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 16645.445
$ws.Range("I62").Value = 9180
$ws.Range("J62").Value = 25977.25
$ws.Range("K62").Value = 9180
$ws.Range("L62").Value = 25977.25
$ws.Range("M62").Value = -8556
$ws.Range("N62").Value = -27225.25
$ws.Range("H65").Value = 16645.445
$ws.Range("I65").Value = 9180
$ws.Range("J65").Value = 25977.25
$ws.Range("K65").Value = 45900
$ws.Range("L65").Value = 129886.25
$ws.Range("M65").Value = -42780
$ws.Range("N65").Value = -136126.25
$ws.Range("H105").Value = 2101
$ws.Range("I105").Value = 2125.6667
$ws.Range("J105").Value = 1990
$ws.Range("K105").Value = 2125.6667
$ws.Range("L105").Value = 1990
$ws.Range("M105").Value = -378.6667000000002
$ws.Range("N105").Value = -5484
$ws.Range("H122").Value = 11038.64
$ws.Range("I122").Value = 973.1667
$ws.Range("J122").Value = 36921.285
$ws.Range("K122").Value = 2919.5001
$ws.Range("L122").Value = 110763.855
$ws.Range("M122").Value = -469.5001000000002
$ws.Range("N122").Value = -115663.855
$ws.Range("H132").Value = 2357.8572
$ws.Range("I132").Value = 2364.1
$ws.Range("J132").Value = 2342.25
$ws.Range("K132").Value = 7092.299999999999
$ws.Range("L132").Value = 7026.75
$ws.Range("M132").Value = -4562.299999999999
$ws.Range("N132").Value = -12086.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 18331.25
$ws.Range("J32").Value = 18331.25
$ws.Range("L32").Value = 18331.25
$ws.Range("N32").Value = -18923.25
$ws.Range("H42").Value = 24632
$ws.Range("J42").Value = 24632
$ws.Range("L42").Value = 24632
$ws.Range("N42").Value = -25602
$ws.Range("H115").Value = 24632
$ws.Range("J115").Value = 24632
$ws.Range("L115").Value = 24632
$ws.Range("N115").Value = -26982
$ws.Range("H122").Value = 4089.7144
$ws.Range("I122").Value = 4303.4165
$ws.Range("J122").Value = 2807.5
$ws.Range("K122").Value = 12910.2495
$ws.Range("L122").Value = 8422.5
$ws.Range("M122").Value = -10460.2495
$ws.Range("N122").Value = -13322.5
$ws.Range("H126").Value = 16749
$ws.Range("I126").Value = 29999
$ws.Range("J126").Value = 3499
$ws.Range("K126").Value = 89997
$ws.Range("L126").Value = 10497
$ws.Range("M126").Value = -87527
$ws.Range("N126").Value = -15437

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6004.1177
$ws.Range("I7").Value = 6957.8
$ws.Range("J7").Value = 4641.7144
$ws.Range("K7").Value = 6957.8
$ws.Range("L7").Value = 4641.7144
$ws.Range("M7").Value = -6845.8
$ws.Range("N7").Value = -4865.7144
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H31").Value = 6485.75
$ws.Range("I31").Value = 974.25
$ws.Range("J31").Value = 11997.25
$ws.Range("K31").Value = 974.25
$ws.Range("L31").Value = 11997.25
$ws.Range("M31").Value = -726.25
$ws.Range("N31").Value = -12493.25
$ws.Range("H46").Value = 1817.95
$ws.Range("I46").Value = 1100
$ws.Range("K46").Value = 1100
$ws.Range("M46").Value = -912
$ws.Range("H50").Value = 24318
$ws.Range("I50").Value = 11020
$ws.Range("J50").Value = 28750.666
$ws.Range("K50").Value = 11020
$ws.Range("L50").Value = 28750.666
$ws.Range("M50").Value = -10383
$ws.Range("N50").Value = -30024.666
$ws.Range("H61").Value = 2243.1667
$ws.Range("I61").Value = 1116.375
$ws.Range("J61").Value = 4496.75
$ws.Range("K61").Value = 1116.375
$ws.Range("L61").Value = 4496.75
$ws.Range("M61").Value = -914.375
$ws.Range("N61").Value = -4900.75
$ws.Range("H68").Value = 14335.333
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251
$ws.Range("H71").Value = 14335.333
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256
$ws.Range("H93").Value = 1557.7646
$ws.Range("I93").Value = 823.5
$ws.Range("J93").Value = 1783.6923
$ws.Range("K93").Value = 823.5
$ws.Range("L93").Value = 1783.6923
$ws.Range("M93").Value = 424.5
$ws.Range("N93").Value = -4279.6923
$ws.Range("H113").Value = 2243.1667
$ws.Range("I113").Value = 1116.375
$ws.Range("J113").Value = 4496.75
$ws.Range("K113").Value = 1116.375
$ws.Range("L113").Value = 4496.75
$ws.Range("M113").Value = 1053.625
$ws.Range("N113").Value = -8836.75
$ws.Range("H122").Value = 3594
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3594
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 10782
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -15682
$ws.Range("H126").Value = 6004.1177
$ws.Range("I126").Value = 6957.8
$ws.Range("J126").Value = 4641.7144
$ws.Range("K126").Value = 20873.4
$ws.Range("L126").Value = 13925.1432
$ws.Range("M126").Value = -18403.4
$ws.Range("N126").Value = -18865.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 24431.5
$ws.Range("J115").Value = 24431.5
$ws.Range("L115").Value = 24431.5
$ws.Range("N115").Value = -27565.5
$ws.Range("H116").Value = 69990
$ws.Range("J116").Value = 69990
$ws.Range("L116").Value = 69990
$ws.Range("N116").Value = -79168
$ws.Range("H122").Value = 2804.92
$ws.Range("I122").Value = 2970.8125
$ws.Range("J122").Value = 2510
$ws.Range("K122").Value = 8912.4375
$ws.Range("L122").Value = 7530
$ws.Range("M122").Value = -6462.4375
$ws.Range("N122").Value = -12430
